# Auto-generated edit script applying the Tonberry_Profits market-data refresh diff.
# For each touched cell we set the literal new value (data comes from an external
# market-board fetch, not Excel formulas -- no <f> formula cells exist in this workbook).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 15465.934
$ws.Range("I18").Value = 13500
$ws.Range("J18").Value = 16776.555
$ws.Range("K18").Value = 13500
$ws.Range("L18").Value = 16776.555
$ws.Range("M18").Value = -13216
$ws.Range("N18").Value = -17344.555
$ws.Range("H51").Value = 6211.4287
$ws.Range("J51").Value = 5000
$ws.Range("L51").Value = 5000
$ws.Range("N51").Value = -5968
$ws.Range("H55").Value = 363.66666
$ws.Range("I55").Value = 295
$ws.Range("J55").Value = 501
$ws.Range("K55").Value = 295
$ws.Range("L55").Value = 501
$ws.Range("M55").Value = -81
$ws.Range("N55").Value = -929
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H113").Value = 19532.334
$ws.Range("I113").Value = 22798.8
$ws.Range("K113").Value = 22798.8
$ws.Range("M113").Value = -19544.8
$ws.Range("H131").Value = 2314.8
$ws.Range("I131").Value = 732.2222
$ws.Range("K131").Value = 2196.6666
$ws.Range("M131").Value = 2843.3334
$ws.Range("H137").Value = 1768.3572
$ws.Range("I137").Value = 1480.4
$ws.Range("K137").Value = 4441.200000000001
$ws.Range("M137").Value = -1891.200000000001
$ws.Range("H141").Value = 3223.5454
$ws.Range("I141").Value = 2307.375
$ws.Range("K141").Value = 6922.125
$ws.Range("M141").Value = -1742.125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3083.9553
$ws.Range("I32").Value = 1776.9056
$ws.Range("K32").Value = 1776.9056
$ws.Range("M32").Value = -1489.9056
$ws.Range("H97").Value = 529.5714
$ws.Range("I97").Value = 529.5714
$ws.Range("K97").Value = 529.5714
$ws.Range("M97").Value = -33.57140000000004
$ws.Range("H122").Value = 1374.1316
$ws.Range("I122").Value = 1314.5938
$ws.Range("J122").Value = 1691.6666
$ws.Range("K122").Value = 3943.7814
$ws.Range("L122").Value = 5074.9998
$ws.Range("M122").Value = -1493.7814
$ws.Range("N122").Value = -9974.9998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1000.46155
$ws.Range("I94").Value = 849.625
$ws.Range("J94").Value = 1241.8
$ws.Range("K94").Value = 849.625
$ws.Range("L94").Value = 1241.8
$ws.Range("M94").Value = -398.625
$ws.Range("N94").Value = -2143.8
$ws.Range("H99").Value = 1587.5
$ws.Range("J99").Value = 2000
$ws.Range("L99").Value = 2000
$ws.Range("N99").Value = -4996
$ws.Range("H105").Value = 2387.5
$ws.Range("I105").Value = 2387.5
$ws.Range("K105").Value = 2387.5
$ws.Range("M105").Value = -640.5
$ws.Range("H107").Value = 2255.1428
$ws.Range("I107").Value = 1987.125
$ws.Range("J107").Value = 2612.5
$ws.Range("K107").Value = 1987.125
$ws.Range("L107").Value = 2612.5
$ws.Range("M107").Value = -67.125
$ws.Range("N107").Value = -6452.5
$ws.Range("H132").Value = 99995
$ws.Range("J132").Value = 99995
$ws.Range("L132").Value = 99995
$ws.Range("N132").Value = -110115
$ws.Range("H140").Value = 49999.668
$ws.Range("J140").Value = 49999.668
$ws.Range("L140").Value = 49999.668
$ws.Range("N140").Value = -60359.668

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 665.4
$ws.Range("I105").Value = 665.4
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 665.4
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 1081.6
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 466.8095
$ws.Range("I107").Value = 396.7857
$ws.Range("K107").Value = 396.7857
$ws.Range("M107").Value = 1523.2143
$ws.Range("H141").Value = 24481.4
$ws.Range("I141").Value = 15000
$ws.Range("J141").Value = 25534.889
$ws.Range("K141").Value = 15000
$ws.Range("L141").Value = 25534.889
$ws.Range("M141").Value = -9820
$ws.Range("N141").Value = -35894.889

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 586
$ws.Range("J5").Value = 934.6667
$ws.Range("L5").Value = 2804.0001
$ws.Range("N5").Value = -3028.0001
$ws.Range("H131").Value = 12520650
$ws.Range("J131").Value = 25679.406
$ws.Range("L131").Value = 77038.21799999999
$ws.Range("N131").Value = -87118.21799999999
$ws.Range("H132").Value = 1525.875
$ws.Range("J132").Value = 1469.5
$ws.Range("L132").Value = 13225.5
$ws.Range("N132").Value = -18285.5
$ws.Range("H135").Value = 586
$ws.Range("J135").Value = 934.6667
$ws.Range("L135").Value = 8412.0003
$ws.Range("N135").Value = -13482.0003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 150
$ws.Range("I2").Value = 184.83333
$ws.Range("J2").Value = 97.75
$ws.Range("K2").Value = 184.83333
$ws.Range("L2").Value = 97.75
$ws.Range("M2").Value = -71.83332999999999
$ws.Range("N2").Value = -323.75
$ws.Range("H70").Value = 5286.5713
$ws.Range("I70").Value = 5876.5
$ws.Range("K70").Value = 5876.5
$ws.Range("M70").Value = -5606.5
$ws.Range("H73").Value = 5286.5713
$ws.Range("I73").Value = 5876.5
$ws.Range("K73").Value = 5876.5
$ws.Range("M73").Value = -4940.5
$ws.Range("H107").Value = 879.44446
$ws.Range("J107").Value = 1118.8334
$ws.Range("L107").Value = 1118.8334
$ws.Range("N107").Value = -4958.8334
$ws.Range("H122").Value = 1902.7084
$ws.Range("I122").Value = 1770.4667
$ws.Range("J122").Value = 2123.111
$ws.Range("K122").Value = 5311.4001
$ws.Range("L122").Value = 6369.333
$ws.Range("M122").Value = -2861.4001
$ws.Range("N122").Value = -11269.333

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5392.706
$ws.Range("I7").Value = 2542.1428
$ws.Range("J7").Value = 7388.1
$ws.Range("K7").Value = 2542.1428
$ws.Range("L7").Value = 7388.1
$ws.Range("M7").Value = -2430.1428
$ws.Range("N7").Value = -7612.1
$ws.Range("H40").Value = 8592
$ws.Range("I40").Value = 5761.8335
$ws.Range("J40").Value = 10135.728
$ws.Range("K40").Value = 5761.8335
$ws.Range("L40").Value = 10135.728
$ws.Range("M40").Value = -5625.8335
$ws.Range("N40").Value = -10407.728
$ws.Range("H46").Value = 1773.0769
$ws.Range("J46").Value = 1773.0769
$ws.Range("L46").Value = 1773.0769
$ws.Range("N46").Value = -2149.0769
$ws.Range("H55").Value = 309.9375
$ws.Range("I55").Value = 112.71429
$ws.Range("J55").Value = 463.33334
$ws.Range("K55").Value = 112.71429
$ws.Range("L55").Value = 463.33334
$ws.Range("M55").Value = 60.28570999999999
$ws.Range("N55").Value = -809.33334
$ws.Range("H122").Value = 6282.28
$ws.Range("J122").Value = 7643.2144
$ws.Range("L122").Value = 22929.6432
$ws.Range("N122").Value = -27829.6432
$ws.Range("H126").Value = 5392.706
$ws.Range("I126").Value = 2542.1428
$ws.Range("J126").Value = 7388.1
$ws.Range("K126").Value = 7626.428400000001
$ws.Range("L126").Value = 22164.3
$ws.Range("M126").Value = -5156.428400000001
$ws.Range("N126").Value = -27104.3

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 22708.916
$ws.Range("I122").Value = 33020.125
$ws.Range("J122").Value = 2086.5
$ws.Range("K122").Value = 99060.375
$ws.Range("L122").Value = 6259.5
$ws.Range("M122").Value = -96610.375
$ws.Range("N122").Value = -11159.5
$ws.Range("H125").Value = 39642.855
$ws.Range("J125").Value = 39642.855
$ws.Range("L125").Value = 39642.855
$ws.Range("N125").Value = -49482.855
$ws.Range("H126").Value = 3777.743
$ws.Range("I126").Value = 2928.68
$ws.Range("J126").Value = 5900.4
$ws.Range("K126").Value = 8786.039999999999
$ws.Range("L126").Value = 17701.2
$ws.Range("M126").Value = -6316.039999999999

Write-Output "Applied 196 cell edits across 8 sheets."
